# Beta version of DataParser
#
# - Rename the "R1_OUT" entry on the Source sheet to "R1out".
# - Make the Source sheet the active tab/sheet (previously EB was active),
#   with cell A2 selected (previously B3 was selected on Source).

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Source")
$source.Range("A2").Value = "R1out"

$source.Activate()
$source.Range("A2").Select()
